$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.967.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.04%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.604"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.88%  "

$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.990"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.651.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.299.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.859.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.12%  "

$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("E23").Value = "  -2.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("E31").Value = "  -1.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0866"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.96%  "

$ws.Range("E34").Value = "  +6.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.02%  "

$ws.Range("E36").Value = "  -1.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0356"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.64%  "

$ws.Range("E40").Value = "  -3.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("E44").Value = "  +0.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.708.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.96%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.49%  "
